$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 529, shifting existing rows 529:595 down to 530:596.
$ws.Rows.Item(529).Insert()

# Populate the newly inserted row 529 with the new record's data.
$ws.Range("A529").Value = 3
$ws.Range("B529").Value = "Femacal de La Calera"
$ws.Range("C529").Value = "Coquimbo"
$ws.Range("D529").Value = 45142
$ws.Range("E529").Value = 5
$ws.Range("F529").Value = 100112009
$ws.Range("G529").Value = "Acelga"
$ws.Range("H529").Value = "Sin especificar"
$ws.Range("I529").Value = "Primera"
$ws.Range("J529").Value = 140
$ws.Range("K529").Value = 3000
$ws.Range("L529").Value = 3000
$ws.Range("M529").Value = 3000
$ws.Range("N529").Value = "$/docena de atados (6 kilos)"
$ws.Range("O529").Value = "Provincia de Quillota"
$ws.Range("P529").Value = 500
$ws.Range("Q529").Value = 6
$ws.Range("R529").Value = "Hortaliza"
